# Add 2022-Q4 data
#
# 1) "总计" (Total) sheet: insert a new row right after the header for the
#    2022-Q4 summary figures, shifting the existing quarters down by one
#    row, and renumber the leading index column.
# 2) Insert a brand-new "2022-Q4" worksheet right after "总计" holding the
#    per-fund detail for that quarter (cloned from an existing quarter
#    sheet so the styles/number formats line up, then its values are
#    overwritten).

$wb = $excel.ActiveWorkbook

# --- 1) Update the "总计" (总计) summary sheet -----------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.02

for ($r = 3; $r -le 10; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# --- 2) Insert the new "2022-Q4" detail sheet ------------------------------
# Clone "2022-Q3" (current position 2) so the header row / cell styles match
# the other quarter sheets, then drop it right after "总计".
$template = $wb.Worksheets.Item(2)
$template.Copy($null, $wb.Worksheets.Item(1))

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template has two data rows; 2022-Q4 only needs one.
$q4.Rows.Item(3).Delete()

# Force the fund-code/amount columns back to text (matches the source data,
# which stores these as strings, e.g. so leading zeros survive) before
# writing the new values.
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "004351"
$q4.Cells.Item(2, 3).Value = "汇丰晋信珠三角区域发展混合"
$q4.Cells.Item(2, 4).Value = "0.52"
$q4.Cells.Item(2, 5).Value = "93.91"
$q4.Cells.Item(2, 6).Value = "3.40"
$q4.Cells.Item(2, 7).Value = "0.0177"
$q4.Cells.Item(2, 8).Value = 6
